$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump "Latest HO Xliff Generate Date" for rows 2 & 3 ---
# (they share the same underlying value, so both rows move together)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 10:16:28"
$wsOverview.Range("G3").Value = "2016-08-23 10:16:28"

# --- zh-cn sheet: Priority ht -> mt, refreshed handoff/handback datetimes ---
# rows 2 & 3 shared the same values, so both rows move together
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-23 10:16:23"
$wsZhCn.Range("K2").Value = "2016-08-23 10:16:40"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-23 10:16:23"
$wsZhCn.Range("K3").Value = "2016-08-23 10:16:40"

# --- de-de sheet: Priority ht -> mt, refreshed handoff/handback datetimes ---
# Correspond Handoff Datetime here shares its value with the Overview sheet's
# Latest HO Xliff Generate Date, and Correspond Handback DateTime is refreshed too.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-23 10:16:28"
$wsDeDe.Range("K2").Value = "2016-08-23 10:16:46"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-23 10:16:28"
$wsDeDe.Range("K3").Value = "2016-08-23 10:16:46"
